$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "email" column header (E1)
$ws.Range("E1").Value = "email"

# Email address per row (matched by first name), written first so the
# hyperlink keeps this as its display text, then turned into a mailto:
# hyperlink - this mints the "Hyperlink" cell style (underlined,
# theme-colored font) used by Excel for hyperlinked cells.
$ws.Range("E2").Value = "harry@gmail.com"
$ws.Range("E3").Value = "jane@gmail.com"
$ws.Range("E4").Value = "alex@gmail.com"
$ws.Range("E5").Value = "chris@gmail.com"
$ws.Range("E6").Value = "peter@gmail.com"

$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:harry@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:jane@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:alex@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:chris@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:peter@gmail.com")

# Widen the new column to fit the email addresses
$ws.Columns("E").ColumnWidth = 21.8

# Reflect the cursor resting on F12 after the edits (as in the saved file)
$ws.Range("F12").Select()
